$wb = $excel.ActiveWorkbook

# Citywide Totals
$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("J2").Value = 7144
$ws.Range("J3").Value = 7541
$ws.Range("J4").Value = 1645
$ws.Range("J5").Value = 588
$ws.Range("J6").Value = 10256
$ws.Range("J7").Value = 27174

# By Neighborhood
$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J4").Value = 124
$ws.Range("J6").Value = 210
$ws.Range("J7").Value = 776
$ws.Range("J8").Value = 1707
$ws.Range("J9").Value = 140
$ws.Range("J11").Value = 487
$ws.Range("J15").Value = 338
$ws.Range("J20").Value = 572
$ws.Range("J23").Value = 247
$ws.Range("J27").Value = 167
$ws.Range("J29").Value = 1452
$ws.Range("J31").Value = 283
$ws.Range("J33").Value = 1234
$ws.Range("J37").Value = 834
$ws.Range("J42").Value = 1163
$ws.Range("J43").Value = 229
$ws.Range("J44").Value = 212
$ws.Range("J46").Value = 92
$ws.Range("J47").Value = 200
$ws.Range("J49").Value = 166
$ws.Range("J50").Value = 162
$ws.Range("J52").Value = 692
$ws.Range("J53").Value = 400
$ws.Range("J54").Value = 534
$ws.Range("J63").Value = 83
$ws.Range("J65").Value = 679
$ws.Range("J67").Value = 1013
$ws.Range("J77").Value = 189
$ws.Range("J78").Value = 314
$ws.Range("J79").Value = 749
$ws.Range("J83").Value = 545
$ws.Range("J84").Value = 224
$ws.Range("J85").Value = 1117
$ws.Range("J88").Value = 289
$ws.Range("J90").Value = 285
$ws.Range("J91").Value = 311
$ws.Range("J92").Value = 90
$ws.Range("J94").Value = 297
$ws.Range("J95").Value = 390
$ws.Range("J98").Value = 202
$ws.Range("J99").Value = 419
$ws.Range("J101").Value = 27174

# Auburn Gresham
$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("J2").Value = 242
$ws.Range("J6").Value = 248
$ws.Range("J7").Value = 776

# Belmont Cragin
$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("J4").Value = 28
$ws.Range("J6").Value = 230
$ws.Range("J7").Value = 487

# South Shore
$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("J2").Value = 299
$ws.Range("J3").Value = 404
$ws.Range("J7").Value = 1117

# Little Village
$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("J6").Value = 297
$ws.Range("J7").Value = 692

# Logan Square
$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("J3").Value = 47
$ws.Range("J6").Value = 266
$ws.Range("J7").Value = 400

# Austin
$ws = $wb.Worksheets.Item("Austin")
$ws.Range("J2").Value = 450
$ws.Range("J6").Value = 624
$ws.Range("J7").Value = 1707

# South Chicago
$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("J2").Value = 159
$ws.Range("J7").Value = 545

# Garfield Park
$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("J2").Value = 277
$ws.Range("J3").Value = 409
$ws.Range("J6").Value = 441
$ws.Range("J7").Value = 1234

# West Pullman
$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("J2").Value = 137
$ws.Range("J3").Value = 140
$ws.Range("J7").Value = 390

# Grand Crossing
$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("J2").Value = 250
$ws.Range("J6").Value = 241
$ws.Range("J7").Value = 834

# New City
$ws = $wb.Worksheets.Item("New City")
$ws.Range("J2").Value = 195
$ws.Range("J6").Value = 252
$ws.Range("J7").Value = 679

# Woodlawn
$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("J2").Value = 114
$ws.Range("J6").Value = 110
$ws.Range("J7").Value = 419

# Gage Park
$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("J3").Value = 67
$ws.Range("J6").Value = 101
$ws.Range("J7").Value = 283

# North Lawndale
$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("J2").Value = 258
$ws.Range("J3").Value = 379
$ws.Range("J7").Value = 1013

# South Deering
$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("J6").Value = 74
$ws.Range("J7").Value = 224

# Lincoln Park
$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("J6").Value = 94
$ws.Range("J7").Value = 166

# Loop
$ws = $wb.Worksheets.Item("Loop")
$ws.Range("J2").Value = 134
$ws.Range("J3").Value = 107
$ws.Range("J6").Value = 246
$ws.Range("J7").Value = 534

# Englewood
$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("J2").Value = 441
$ws.Range("J3").Value = 513
$ws.Range("J4").Value = 78
$ws.Range("J5").Value = 54
$ws.Range("J7").Value = 1452

# Irving Park
$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("J2").Value = 65
$ws.Range("J3").Value = 48
$ws.Range("J6").Value = 84
$ws.Range("J7").Value = 212

# Ashburn
$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("J4").Value = 15
$ws.Range("J7").Value = 210

# Humboldt Park
$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("J3").Value = 232
$ws.Range("J6").Value = 620
$ws.Range("J7").Value = 1163

# Rogers Park
$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("J2").Value = 84
$ws.Range("J6").Value = 96
$ws.Range("J7").Value = 314

# Jefferson Park
$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Range("J2").Value = 25
$ws.Range("J6").Value = 38
$ws.Range("J7").Value = 92

# Douglas
$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("J6").Value = 68
$ws.Range("J7").Value = 247

# Washington Park
$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("J2").Value = 82
$ws.Range("J3").Value = 128
$ws.Range("J7").Value = 311

# Roseland
$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("J3").Value = 251
$ws.Range("J7").Value = 749

# Chicago Lawn
$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("J6").Value = 165
$ws.Range("J7").Value = 572

# West Loop
$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("J3").Value = 57
$ws.Range("J7").Value = 297

# Kenwood
$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("J2").Value = 45
$ws.Range("J7").Value = 200

# Brighton Park
$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("J2").Value = 95
$ws.Range("J6").Value = 155
$ws.Range("J7").Value = 338

# Wicker Park
$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("J3").Value = 27
$ws.Range("J7").Value = 202

# Lincoln Square
$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("J3").Value = 41
$ws.Range("J7").Value = 162

# Avalon Park
$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("J3").Value = 45
$ws.Range("J7").Value = 140

# West Elsdon
$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("J6").Value = 31
$ws.Range("J7").Value = 90

# United Center
$ws = $wb.Worksheets.Item("United Center")
$ws.Range("J6").Value = 150
$ws.Range("J7").Value = 289

# Edgewater
$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("J6").Value = 60
$ws.Range("J7").Value = 167

# Washington Heights
$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("J3").Value = 78
$ws.Range("J7").Value = 285

# Hyde Park
$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("J2").Value = 26
$ws.Range("J7").Value = 229

# Riverdale
$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("J3").Value = 62
$ws.Range("J7").Value = 189

# Archer Heights
$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("J3").Value = 25
$ws.Range("J7").Value = 124
